# Regenerate orders with updated distance/sizes.
# The workbook stores trial-order data where text values encode a
# Distance code (D64/D80/D51) and a Size code (S20/S25/S30) as part of
# Condition / Filename_Left / Filename_Right / Distance / Size columns.
# This script renumbers the distance/size codes:
#   D64 -> D69
#   D80 -> D86
#   D51 -> D55
#   S30 -> S31
# (S25 and S20 are left untouched.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange

# Order matters only in that none of the replacement tokens collide with
# any of the search tokens, so a straightforward sequential replace is
# safe and will not double-substitute.
$replacements = @(
    @("D64", "D69"),
    @("D80", "D86"),
    @("D51", "D55"),
    @("S30", "S31")
)

foreach ($pair in $replacements) {
    $find = $pair[0]
    $replace = $pair[1]
    # xlWhole? No - these tokens appear as substrings inside larger
    # strings (e.g. Face07_D80_S25, Fixation_D80_l.png), so use a part
    # match (xlPart = 2) rather than whole-cell match (xlWhole = 1).
    $used.Replace($find, $replace, 2, 1, $false, $false, $false, $false) | Out-Null
}
